# Generate Report for Archive
#
# The localization status report is regenerated: the file
# "6c7127f0-f445-4919-8766-5b113b957ca0.md" moved from "Ready for handoff"
# into "In Translation" (its handoff timestamp reset to the not-yet-handed-off
# placeholder on the Overview-less detail sheets the other file's values show),
# and, because the report rows are keyed off status/order, the two rows
# previously belonging to "064e0f25-b0eb-4c64-af88-7a51c01369bc.md" and
# "6c7127f0-f445-4919-8766-5b113b957ca0.md" swap places on every sheet.
#
# This script reproduces that by writing the new cell values directly
# (values only - no structural/style changes), matching the target OOXML.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": File Name / zh-cn / de-de status columns
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "6c7127f0-f445-4919-8766-5b113b957ca0.md"
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"

$wsOverview.Range("A5").Value = "064e0f25-b0eb-4c64-af88-7a51c01369bc.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet "zh-cn": detail rows for the same two files
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "6c7127f0-f445-4919-8766-5b113b957ca0.md"
$wsZhCn.Range("B4").Value = "In Translation"
$wsZhCn.Range("C4").Value = "6c7127f0-f445-4919-8766-5b113b957ca0.64b513a3b303b40b6cd81761a6b4a1281cec63aa.zh-cn.xlf"
$wsZhCn.Range("D4").Value = "2016-03-08 06:19:05"

$wsZhCn.Range("A5").Value = "064e0f25-b0eb-4c64-af88-7a51c01369bc.md"
$wsZhCn.Range("B5").Value = "Ready for handoff"
$wsZhCn.Range("C5").Value = "064e0f25-b0eb-4c64-af88-7a51c01369bc.1b16b041d8163f91ba6766c3983dac0eba27343c.zh-cn.xlf"
$wsZhCn.Range("D5").Value = "2016-03-08 06:18:12"

# ---------------------------------------------------------------------
# Sheet "de-de": detail rows for the same two files
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "6c7127f0-f445-4919-8766-5b113b957ca0.md"
$wsDeDe.Range("B4").Value = "In Translation"
$wsDeDe.Range("C4").Value = "6c7127f0-f445-4919-8766-5b113b957ca0.64b513a3b303b40b6cd81761a6b4a1281cec63aa.de-de.xlf"
$wsDeDe.Range("D4").Value = "2016-03-08 06:19:07"

$wsDeDe.Range("A5").Value = "064e0f25-b0eb-4c64-af88-7a51c01369bc.md"
$wsDeDe.Range("B5").Value = "Ready for handoff"
$wsDeDe.Range("C5").Value = "064e0f25-b0eb-4c64-af88-7a51c01369bc.1b16b041d8163f91ba6766c3983dac0eba27343c.de-de.xlf"
$wsDeDe.Range("D5").Value = "2016-03-08 06:18:15"
